$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(254).Insert()

$ws.Cells.Item(254, 1).Value = 9
$ws.Cells.Item(254, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(254, 3).Value = "Metropolitana"
$ws.Cells.Item(254, 4).Value = 45141
$ws.Cells.Item(254, 5).Value = 13
$ws.Cells.Item(254, 6).Value = 100112026
$ws.Cells.Item(254, 7).Value = "Haba"
$ws.Cells.Item(254, 8).Value = "Sin especificar"
$ws.Cells.Item(254, 9).Value = "Primera"
$ws.Cells.Item(254, 10).Value = 70
$ws.Cells.Item(254, 11).Value = 19000
$ws.Cells.Item(254, 12).Value = 20000
$ws.Cells.Item(254, 13).Value = 19500
$ws.Cells.Item(254, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(254, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(254, 16).Value = 780
$ws.Cells.Item(254, 17).Value = 25
$ws.Cells.Item(254, 18).Value = "Hortaliza"
